$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 0.407
$ws.Range("L3").Value = 0.447
$ws.Range("C4").Value = 0.68
$ws.Range("I5").Value = 0.551
$ws.Range("K5").Value = 0.417
$ws.Range("E6").Value = 0.598
$ws.Range("H8").Value = 0.404
$ws.Range("E13").Value = 0.613
$ws.Range("G13").Value = 0.759
$ws.Range("B49").Value = 0.688
$ws.Range("G62").Value = 0.543
$ws.Range("E65").Value = 0.603
$ws.Range("B68").Value = 0.634
$ws.Range("K68").Value = 0.469
$ws.Range("F81").Value = 0.568
$ws.Range("G81").Value = 0.69
$ws.Range("H84").Value = 0.423
$ws.Range("C88").Value = 0.562
$ws.Range("L89").Value = 0.398
$ws.Range("F93").Value = 0.498
$ws.Range("G93").Value = 0.595
$ws.Range("H97").Value = 0.368
$ws.Range("J97").Value = 0.339
$ws.Range("L100").Value = 0.401
